# Adds the "Nhung van de khi tu dong hoa IE bang VBScript/Windows Script"
# section to Sheet1: a lead-in sentence in B43, seven bullet points in
# B45:B48 and B49:B52, and a closing summary paragraph in B54.
# This appends 10 new shared strings (indices 29-38) and grows the used
# range from A2:D39 to A2:D54, matching the authored worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newParagraphs = @(
    "Khi tự động hóa trình duyệt Internet Explorer (IE) bằng VBScript và Windows Script, bạn có thể gặp phải một số vấn đề sau:",
    "Khả năng tương thích: IE không còn là trình duyệt phát triển chính thống, và không hỗ trợ nhiều tính năng và tiêu chuẩn web hiện đại. Điều này có thể làm cho các kịch bản tự động hóa trở nên khó khăn vì trang web có thể không hoạt động đúng cách trên IE.",
    "Bảo mật: IE đã gặp nhiều vấn đề bảo mật trong quá khứ và vẫn có thể bị mục tiêu của các cuộc tấn công. Sử dụng IE trong kịch bản tự động hóa có thể đặt dự án của bạn vào nguy cơ bảo mật.",
    "Khả năng thất bại: Các thay đổi không kiểm soát trong trình duyệt, cập nhật bất ngờ hoặc tương tác không mong muốn có thể dẫn đến việc mã tự động hóa không hoạt động.",
    "Khả năng gỡ lỗi: Gỡ lỗi VBScript và Windows Script có thể khó khăn hơn so với các ngôn ngữ lập trình hiện đại khác, đặc biệt khi bạn cần theo dõi các vấn đề trong tương tác với trình duyệt.",
    "Hiệu suất kém: IE thường chạy chậm hơn so với các trình duyệt hiện đại khác, điều này có thể làm cho tự động hóa chậm và không hiệu quả.",
    "Khả năng mở rộng hạn chế: VBScript không phải là một ngôn ngữ lập trình mạnh mẽ như Python hoặc JavaScript, điều này có thể giới hạn khả năng xây dựng các kịch bản tự động hóa phức tạp.",
    "Khả năng tương thích với các trình duyệt khác: VBScript thường được thiết kế cho IE và có thể không hoạt động tốt với các trình duyệt khác như Chrome, Firefox hoặc Edge.",
    "Khả năng tương thích đa nền tảng: VBScript và Windows Script thường chỉ hỗ trợ trên Windows, do đó, không phù hợp cho các dự án đa nền tảng.",
    "Tóm lại, việc tự động hóa trình duyệt IE bằng VBScript và Windows Script có thể đối mặt với nhiều vấn đề về tương thích, hiệu suất và bảo mật. Để thực hiện tự động hóa hiệu quả hơn, bạn nên xem xét sử dụng các công cụ hiện đại hơn như Selenium WebDriver và các ngôn ngữ lập trình như Python, JavaScript hoặc C#."
)

# Target rows in column B for each paragraph, in order (row 44 and 53 stay blank).
$targetRows = @(43, 45, 46, 47, 48, 49, 50, 51, 52, 54)

for ($i = 0; $i -lt $targetRows.Length; $i++) {
    $ws.Cells.Item($targetRows[$i], 2).Value = $newParagraphs[$i]
}

# Match the author's final view/selection state in the saved file.
$ws.Activate()
$ws.Range("AC20").Select()
